$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1643356643356643
$ws.Range("C2").Value = 0.6293706293706294
$ws.Range("J2").Value = 0.02797202797202797
$ws.Range("P2").Value = 0.1223776223776224
$ws.Range("S2").Value = 0.05594405594405594
$ws.Range("B3").Value = 0.005319148936170213
$ws.Range("C3").Value = 0.01595744680851064
$ws.Range("J3").Value = 0.04787234042553191
$ws.Range("P3").Value = 0.8085106382978723
$ws.Range("S3").Value = 0.1223404255319149
$ws.Range("J4").Value = 0.04545454545454546
$ws.Range("P4").Value = 0.7727272727272727
$ws.Range("S4").Value = 0.1818181818181818
$ws.Range("B6").Value = 0.04587155963302753
$ws.Range("D6").Value = 0.02752293577981652
$ws.Range("F6").Value = 0.05045871559633028
$ws.Range("J6").Value = 0.2752293577981652
$ws.Range("O6").Value = 0.01834862385321101
$ws.Range("Q6").Value = 0.1926605504587156
$ws.Range("R6").Value = 0.0871559633027523
$ws.Range("S6").Value = 0.3027522935779817
$ws.Range("B7").Value = 0.07253886010362694
$ws.Range("D7").Value = 0.06217616580310881
$ws.Range("F7").Value = 0.07253886010362694
$ws.Range("J7").Value = 0.1450777202072539
$ws.Range("O7").Value = 0.02590673575129534
$ws.Range("Q7").Value = 0.227979274611399
$ws.Range("R7").Value = 0.07772020725388601
$ws.Range("S7").Value = 0.3160621761658031
$ws.Range("B8").Value = 0.09230769230769231
$ws.Range("D8").Value = 0.01978021978021978
$ws.Range("F8").Value = 0.06813186813186813
$ws.Range("J8").Value = 0.1362637362637363
$ws.Range("O8").Value = 0.02197802197802198
$ws.Range("Q8").Value = 0.1714285714285714
$ws.Range("R8").Value = 0.1032967032967033
$ws.Range("S8").Value = 0.3868131868131868
$ws.Range("B9").Value = 0.08749999999999999
$ws.Range("D9").Value = 0.01875
$ws.Range("F9").Value = 0.08749999999999999
$ws.Range("J9").Value = 0.14375
$ws.Range("O9").Value = 0.01875
$ws.Range("Q9").Value = 0.26875
$ws.Range("R9").Value = 0.05625
$ws.Range("S9").Value = 0.31875
$ws.Range("B10").Value = 0.1064880112834979
$ws.Range("D10").Value = 0.02538787023977433
$ws.Range("F10").Value = 0.06346967559943582
$ws.Range("J10").Value = 0.1248236953455571
$ws.Range("O10").Value = 0.01057827926657264
$ws.Range("Q10").Value = 0.2228490832157969
$ws.Range("R10").Value = 0.09873060648801128
$ws.Range("S10").Value = 0.347672778561354
$ws.Range("G11").Value = 0.1473684210526316
$ws.Range("J11").Value = 0.08771929824561403
$ws.Range("K11").Value = 0.2105263157894737
$ws.Range("L11").Value = 0.543859649122807
$ws.Range("S11").Value = 0.01052631578947368
$ws.Range("G12").Value = 0.7453416149068323
$ws.Range("J12").Value = 0.1677018633540373
$ws.Range("K12").Value = 0.006211180124223602
$ws.Range("L12").Value = 0.04968944099378882
$ws.Range("S12").Value = 0.03105590062111801
$ws.Range("G13").Value = 0.7391304347826086
$ws.Range("J13").Value = 0.2391304347826087
$ws.Range("S13").Value = 0.02173913043478261
$ws.Range("F15").Value = 0.02049180327868852
$ws.Range("H15").Value = 0.1639344262295082
$ws.Range("I15").Value = 0.05737704918032787
$ws.Range("J15").Value = 0.4262295081967213
$ws.Range("K15").Value = 0.05327868852459016
$ws.Range("M15").Value = 0.00819672131147541
$ws.Range("O15").Value = 0.05737704918032787
$ws.Range("S15").Value = 0.2131147540983606
$ws.Range("F16").Value = 0.03167420814479638
$ws.Range("H16").Value = 0.1945701357466063
$ws.Range("I16").Value = 0.05882352941176471
$ws.Range("J16").Value = 0.4434389140271493
$ws.Range("K16").Value = 0.09502262443438914
$ws.Range("M16").Value = 0.02714932126696833
$ws.Range("N16").Value = 0.004524886877828055
$ws.Range("O16").Value = 0.05882352941176471
$ws.Range("S16").Value = 0.08597285067873303
$ws.Range("F17").Value = 0.01529636711281071
$ws.Range("H17").Value = 0.1778202676864245
$ws.Range("I17").Value = 0.08030592734225621
$ws.Range("J17").Value = 0.4416826003824092
$ws.Range("K17").Value = 0.0994263862332696
$ws.Range("M17").Value = 0.01720841300191205
$ws.Range("O17").Value = 0.05927342256214149
$ws.Range("S17").Value = 0.1089866156787763
$ws.Range("F18").Value = 0.01769911504424779
$ws.Range("H18").Value = 0.2035398230088496
$ws.Range("I18").Value = 0.06637168141592921
$ws.Range("J18").Value = 0.4867256637168141
$ws.Range("K18").Value = 0.07964601769911504
$ws.Range("M18").Value = 0.008849557522123894
$ws.Range("O18").Value = 0.05752212389380531
$ws.Range("S18").Value = 0.07964601769911504
$ws.Range("F19").Value = 0.01284246575342466
$ws.Range("H19").Value = 0.2020547945205479
$ws.Range("I19").Value = 0.06335616438356165
$ws.Range("J19").Value = 0.3938356164383562
$ws.Range("K19").Value = 0.0976027397260274
$ws.Range("M19").Value = 0.02568493150684931
$ws.Range("N19").Value = 0.001712328767123288
$ws.Range("O19").Value = 0.08732876712328767
$ws.Range("S19").Value = 0.1155821917808219
